# Update countries & provincias Spain
# Refresh of the COVID-19 "Pais" dashboard: new case counts pulled in for
# several countries, two countries re-ranked (so their name/data swapped
# rows), and the "last updated" timestamp bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 09:52"

# --- Rusia overtakes Austria (rows 19/20 swap identity + new numbers) ---
$ws.Range("A19").Value = "Rusia"
$ws.Range("B19").Value = 15770
$ws.Range("C19").Value = 2186
$ws.Range("D19").Value = 1291
$ws.Range("E19").Value = 14349
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 130

$ws.Range("A20").Value = "Austria"
$ws.Range("B20").Value = 13814
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 6604
$ws.Range("E20").Value = 6873
$ws.Range("F20").Value = 246
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 337

# --- Bosnia y Herzegovina (row 74) refreshed numbers ---
$ws.Range("B74").Value = 975
$ws.Range("C74").Value = 29
$ws.Range("D74").Value = 139
$ws.Range("E74").Value = 799
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 37

# --- Kazajistan (row 75) refreshed numbers ---
$ws.Range("B75").Value = 897
$ws.Range("C75").Value = 32
$ws.Range("D75").Value = 81
$ws.Range("E75").Value = 806
$ws.Range("F75").Value = 21
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 10

# --- Afganistan jumps ahead of Principado de Andorra / Oman / Costa Rica
#     (rows 87-90 shift down one place, Afganistan takes row 87) ---
$ws.Range("A87").Value = "Afganistan"
$ws.Range("B87").Value = 607
$ws.Range("C87").Value = 52
$ws.Range("D87").Value = 32
$ws.Range("E87").Value = 557
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 18

$ws.Range("A88").Value = "Principado de Andorra"
$ws.Range("B88").Value = 601
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 71
$ws.Range("E88").Value = 504
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 26

$ws.Range("A89").Value = "Oman"
$ws.Range("B89").Value = 599
$ws.Range("C89").Value = 53
$ws.Range("D89").Value = 109
$ws.Range("E89").Value = 487
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 3

$ws.Range("A90").Value = "Costa Rica"
$ws.Range("B90").Value = 577
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 49
$ws.Range("E90").Value = 525
$ws.Range("F90").Value = 13
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 3

# --- Montenegro (row 110) refreshed numbers ---
$ws.Range("B110").Value = 267
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 5
$ws.Range("E110").Value = 260
$ws.Range("F110").Value = 7
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

# --- San Bartolome (row 204) refreshed numbers ---
$ws.Range("D204").Value = 4
$ws.Range("E204").Value = 2
